$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws2 = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet: updated forecast numbers (Auto ARIMA removed) ---
$ws1.Cells.Item(2, 3).Value = 55
$ws1.Cells.Item(2, 4).Value = 55
$ws1.Cells.Item(2, 5).Value = 66
$ws1.Cells.Item(2, 6).Value = 77
$ws1.Cells.Item(2, 7).Value = 94

$ws1.Cells.Item(3, 3).Value = 51
$ws1.Cells.Item(3, 4).Value = 52
$ws1.Cells.Item(3, 5).Value = 64
$ws1.Cells.Item(3, 6).Value = 77
$ws1.Cells.Item(3, 7).Value = 97

$ws1.Cells.Item(4, 3).Value = 59
$ws1.Cells.Item(4, 4).Value = 43
$ws1.Cells.Item(4, 5).Value = 52
$ws1.Cells.Item(4, 6).Value = 62
$ws1.Cells.Item(4, 7).Value = 78

$ws1.Cells.Item(5, 3).Value = 63
$ws1.Cells.Item(5, 4).Value = 39
$ws1.Cells.Item(5, 5).Value = 48
$ws1.Cells.Item(5, 6).Value = 57
$ws1.Cells.Item(5, 7).Value = 73

$ws1.Cells.Item(6, 3).Value = 60
$ws1.Cells.Item(6, 4).Value = 39
$ws1.Cells.Item(6, 5).Value = 48
$ws1.Cells.Item(6, 6).Value = 59
$ws1.Cells.Item(6, 7).Value = 77

$ws1.Cells.Item(7, 3).Value = 55
$ws1.Cells.Item(7, 4).Value = 39
$ws1.Cells.Item(7, 5).Value = 48
$ws1.Cells.Item(7, 6).Value = 59
$ws1.Cells.Item(7, 7).Value = 78

$ws1.Cells.Item(8, 3).Value = 54
$ws1.Cells.Item(8, 4).Value = 39
$ws1.Cells.Item(8, 5).Value = 48
$ws1.Cells.Item(8, 6).Value = 60
$ws1.Cells.Item(8, 7).Value = 80

$ws1.Cells.Item(9, 3).Value = 55
$ws1.Cells.Item(9, 4).Value = 39
$ws1.Cells.Item(9, 5).Value = 48
$ws1.Cells.Item(9, 6).Value = 62
$ws1.Cells.Item(9, 7).Value = 85

$ws1.Cells.Item(10, 3).Value = 58
$ws1.Cells.Item(10, 4).Value = 38
$ws1.Cells.Item(10, 5).Value = 46
$ws1.Cells.Item(10, 6).Value = 58
$ws1.Cells.Item(10, 7).Value = 77

$ws1.Cells.Item(11, 3).Value = 67
$ws1.Cells.Item(11, 4).Value = 38
$ws1.Cells.Item(11, 5).Value = 47
$ws1.Cells.Item(11, 6).Value = 61
$ws1.Cells.Item(11, 7).Value = 84

$ws1.Cells.Item(12, 3).Value = 82
$ws1.Cells.Item(12, 4).Value = 39
$ws1.Cells.Item(12, 5).Value = 47
$ws1.Cells.Item(12, 6).Value = 62
$ws1.Cells.Item(12, 7).Value = 85

$ws1.Cells.Item(13, 3).Value = 97
$ws1.Cells.Item(13, 4).Value = 41
$ws1.Cells.Item(13, 5).Value = 50
$ws1.Cells.Item(13, 6).Value = 65
$ws1.Cells.Item(13, 7).Value = 89

$ws1.Cells.Item(14, 3).Value = 98
$ws1.Cells.Item(14, 4).Value = 39
$ws1.Cells.Item(14, 5).Value = 48
$ws1.Cells.Item(14, 6).Value = 62
$ws1.Cells.Item(14, 7).Value = 86

$ws1.Cells.Item(15, 3).Value = 82
$ws1.Cells.Item(15, 4).Value = 38
$ws1.Cells.Item(15, 5).Value = 46
$ws1.Cells.Item(15, 6).Value = 61
$ws1.Cells.Item(15, 7).Value = 86

$ws1.Cells.Item(16, 3).Value = 62
$ws1.Cells.Item(16, 4).Value = 39
$ws1.Cells.Item(16, 5).Value = 48
$ws1.Cells.Item(16, 6).Value = 62
$ws1.Cells.Item(16, 7).Value = 87

$ws1.Cells.Item(17, 3).Value = 55
$ws1.Cells.Item(17, 4).Value = 37
$ws1.Cells.Item(17, 5).Value = 46
$ws1.Cells.Item(17, 6).Value = 60
$ws1.Cells.Item(17, 7).Value = 84

# --- Summary sheet: recomputed aggregate metrics (stored as text, like the originals) ---
$ws2.Cells.Item(9, 2).NumberFormat = "@"
$ws2.Cells.Item(9, 2).Value = "1053"

$ws2.Cells.Item(10, 2).NumberFormat = "@"
$ws2.Cells.Item(10, 2).Value = "452"

$ws2.Cells.Item(11, 2).NumberFormat = "@"
$ws2.Cells.Item(11, 2).Value = "228"

$ws2.Cells.Item(12, 2).NumberFormat = "@"
$ws2.Cells.Item(12, 2).Value = "98"

$ws2.Cells.Item(14, 2).NumberFormat = "@"
$ws2.Cells.Item(14, 2).Value = "51"
